$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 532.0714
$ws.Range("I41").Value = 146
$ws.Range("J41").Value = 1046.8334
$ws.Range("K41").Value = 146
$ws.Range("L41").Value = 1046.8334
$ws.Range("M41").Value = 294
$ws.Range("N41").Value = -1926.8334
$ws.Range("H62").Value = 37271.863
$ws.Range("I62").Value = 63804
$ws.Range("J62").Value = 4616.923
$ws.Range("K62").Value = 63804
$ws.Range("L62").Value = 4616.923
$ws.Range("M62").Value = -63180
$ws.Range("N62").Value = -5864.923
$ws.Range("H65").Value = 37271.863
$ws.Range("I65").Value = 63804
$ws.Range("J65").Value = 4616.923
$ws.Range("K65").Value = 319020
$ws.Range("L65").Value = 23084.615
$ws.Range("M65").Value = -315900
$ws.Range("N65").Value = -29324.615
$ws.Range("H127").Value = 2046.625
$ws.Range("I127").Value = 631
$ws.Range("J127").Value = 2728.2222
$ws.Range("K127").Value = 1893
$ws.Range("L127").Value = 8184.6666
$ws.Range("M127").Value = 3067
$ws.Range("N127").Value = -18104.6666
$ws.Range("H132").Value = 1110.6721
$ws.Range("I132").Value = 1124.8776
$ws.Range("J132").Value = 1052.6666
$ws.Range("K132").Value = 3374.6328
$ws.Range("L132").Value = 3157.9998
$ws.Range("M132").Value = -844.6328000000003
$ws.Range("N132").Value = -8217.9998
$ws.Range("H141").Value = 1791.5714
$ws.Range("I141").Value = 742.0244
$ws.Range("J141").Value = 4660.3335
$ws.Range("K141").Value = 2226.0732
$ws.Range("L141").Value = 13981.0005
$ws.Range("M141").Value = 2953.9268
$ws.Range("N141").Value = -24341.0005

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 16762.766
$ws.Range("I32").Value = 21226.871
$ws.Range("J32").Value = 5602.5
$ws.Range("K32").Value = 21226.871
$ws.Range("L32").Value = 5602.5
$ws.Range("M32").Value = -20939.871
$ws.Range("N32").Value = -6176.5
$ws.Range("H63").Value = 2150.25
$ws.Range("I63").Value = 2093.6
$ws.Range("J63").Value = 3000
$ws.Range("K63").Value = 2093.6
$ws.Range("L63").Value = 3000
$ws.Range("M63").Value = -1407.6
$ws.Range("N63").Value = -4372
$ws.Range("H66").Value = 2150.25
$ws.Range("I66").Value = 2093.6
$ws.Range("J66").Value = 3000
$ws.Range("K66").Value = 10468
$ws.Range("L66").Value = 15000
$ws.Range("M66").Value = -7036
$ws.Range("N66").Value = -21864
$ws.Range("H74").Value = 872.2622699999999
$ws.Range("I74").Value = 787.83673
$ws.Range("K74").Value = 787.83673
$ws.Range("M74").Value = 86.16327000000001
$ws.Range("H77").Value = 872.2622699999999
$ws.Range("I77").Value = 787.83673
$ws.Range("K77").Value = 3939.18365
$ws.Range("M77").Value = 428.8163500000001
$ws.Range("H97").Value = 1051.1111
$ws.Range("I97").Value = 1076.6666
$ws.Range("J97").Value = 1000
$ws.Range("K97").Value = 1076.6666
$ws.Range("L97").Value = 1000
$ws.Range("M97").Value = -580.6666
$ws.Range("N97").Value = -1992
$ws.Range("H122").Value = 968
$ws.Range("I122").Value = 935.3333
$ws.Range("J122").Value = 1033.3334
$ws.Range("K122").Value = 2805.9999
$ws.Range("L122").Value = 3100.0002
$ws.Range("M122").Value = -355.9998999999998
$ws.Range("N122").Value = -8000.0002

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1679.4681
$ws.Range("I86").Value = 1539.6
$ws.Range("J86").Value = 2478.7144
$ws.Range("K86").Value = 1539.6
$ws.Range("L86").Value = 2478.7144
$ws.Range("M86").Value = -416.5999999999999
$ws.Range("N86").Value = -4724.7144
$ws.Range("H89").Value = 1679.4681
$ws.Range("I89").Value = 1539.6
$ws.Range("J89").Value = 2478.7144
$ws.Range("K89").Value = 7698
$ws.Range("L89").Value = 12393.572
$ws.Range("M89").Value = -2082
$ws.Range("N89").Value = -23625.572

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2979.3777
$ws.Range("I31").Value = 3205.68
$ws.Range("J31").Value = 2696.5
$ws.Range("K31").Value = 3205.68
$ws.Range("L31").Value = 2696.5
$ws.Range("M31").Value = -2910.68
$ws.Range("N31").Value = -3286.5
$ws.Range("H34").Value = 2979.3777
$ws.Range("I34").Value = 3205.68
$ws.Range("J34").Value = 2696.5
$ws.Range("K34").Value = 3205.68
$ws.Range("L34").Value = 2696.5
$ws.Range("M34").Value = -3003.68
$ws.Range("N34").Value = -3100.5
$ws.Range("H122").Value = 798
$ws.Range("I122").Value = 500
$ws.Range("K122").Value = 1500
$ws.Range("M122").Value = 950

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 893.62964
$ws.Range("I5").Value = 892.4091
$ws.Range("K5").Value = 2677.2273
$ws.Range("M5").Value = -2565.2273
$ws.Range("H80").Value = 4166.3335
$ws.Range("J80").Value = 4832.8335
$ws.Range("L80").Value = 14498.5005
$ws.Range("N80").Value = -16370.5005
$ws.Range("H83").Value = 4166.3335
$ws.Range("J83").Value = 4832.8335
$ws.Range("L83").Value = 43495.5015
$ws.Range("N83").Value = -52855.5015
$ws.Range("H107").Value = 864521.4399999999
$ws.Range("I107").Value = 189.2
$ws.Range("J107").Value = 1944936.8
$ws.Range("K107").Value = 567.5999999999999
$ws.Range("L107").Value = 5834810.4
$ws.Range("M107").Value = 1352.4
$ws.Range("N107").Value = -5838650.4
$ws.Range("H122").Value = 261.78787
$ws.Range("I122").Value = 177.16667
$ws.Range("J122").Value = 363.33334
$ws.Range("K122").Value = 1594.50003
$ws.Range("L122").Value = 3270.00006
$ws.Range("M122").Value = 855.4999699999998
$ws.Range("N122").Value = -8170.00006
$ws.Range("H131").Value = 6424621
$ws.Range("I131").Value = 100202090
$ws.Range("J131").Value = 1506.8219
$ws.Range("K131").Value = 300606270
$ws.Range("L131").Value = 4520.4657
$ws.Range("M131").Value = -300601230
$ws.Range("N131").Value = -14600.4657
$ws.Range("H135").Value = 893.62964
$ws.Range("I135").Value = 892.4091
$ws.Range("K135").Value = 8031.6819
$ws.Range("M135").Value = -5496.6819

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4399.048
$ws.Range("I70").Value = 4098.1816
$ws.Range("J70").Value = 4730
$ws.Range("K70").Value = 4098.1816
$ws.Range("L70").Value = 4730
$ws.Range("M70").Value = -3828.1816
$ws.Range("N70").Value = -5270
$ws.Range("H73").Value = 4399.048
$ws.Range("I73").Value = 4098.1816
$ws.Range("J73").Value = 4730
$ws.Range("K73").Value = 4098.1816
$ws.Range("L73").Value = 4730
$ws.Range("M73").Value = -3162.1816
$ws.Range("N73").Value = -6602
$ws.Range("H97").Value = 1581.909
$ws.Range("I97").Value = 1218.3334
$ws.Range("J97").Value = 2018.2
$ws.Range("K97").Value = 1218.3334
$ws.Range("L97").Value = 2018.2
$ws.Range("M97").Value = -722.3334
$ws.Range("N97").Value = -3010.2
$ws.Range("H132").Value = 1833.4375
$ws.Range("I132").Value = 1830.0625
$ws.Range("J132").Value = 1840.1875
$ws.Range("K132").Value = 5490.1875
$ws.Range("L132").Value = 5520.5625
$ws.Range("M132").Value = -2960.1875
$ws.Range("N132").Value = -10580.5625

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1457.0769
$ws.Range("I82").Value = 963.1429000000001
$ws.Range("J82").Value = 2033.3334
$ws.Range("K82").Value = 963.1429000000001
$ws.Range("L82").Value = 2033.3334
$ws.Range("M82").Value = -602.1429000000001
$ws.Range("N82").Value = -2755.3334
$ws.Range("H85").Value = 1457.0769
$ws.Range("I85").Value = 963.1429000000001
$ws.Range("J85").Value = 2033.3334
$ws.Range("K85").Value = 963.1429000000001
$ws.Range("L85").Value = 2033.3334
$ws.Range("M85").Value = 284.8570999999999
$ws.Range("N85").Value = -4529.3334
$ws.Range("H122").Value = 7021.3
$ws.Range("I122").Value = 8595.066000000001
$ws.Range("K122").Value = 25785.198
$ws.Range("M122").Value = -23335.198
$ws.Range("H136").Value = 1992.3778
$ws.Range("I136").Value = 1114.1538
$ws.Range("J136").Value = 7700.8335
$ws.Range("K136").Value = 3342.4614
$ws.Range("L136").Value = 23102.5005
$ws.Range("M136").Value = -792.4614000000001
$ws.Range("N136").Value = -28202.5005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 810.75
$ws.Range("I122").Value = 551.5
$ws.Range("J122").Value = 1070
$ws.Range("K122").Value = 1654.5
$ws.Range("L122").Value = 3210
$ws.Range("M122").Value = 795.5
$ws.Range("N122").Value = -8110
$ws.Range("H136").Value = 2050
$ws.Range("I136").Value = 2491.5
$ws.Range("K136").Value = 7474.5
$ws.Range("M136").Value = -4924.5
$ws.Range("H140").Value = 35496.727
$ws.Range("J140").Value = 35496.727
$ws.Range("L140").Value = 35496.727
$ws.Range("N140").Value = -45856.727
$ws.Range("H141").Value = 84590
$ws.Range("J141").Value = 84590
$ws.Range("L141").Value = 84590
$ws.Range("N141").Value = -94950
